$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "socks" paragraph - split the trailing "following:" off into its
# own run wrapped in proofErr gramStart/gramEnd (grammar-check marker around
# a single word), matching Word's "sentence ends oddly" grammar flag.
# ---------------------------------------------------------------------------
$socksPara = $d.Paragraphs.Item(40)
$socksRange = $socksPara.Range
if ($socksRange.Text -notlike "There are 20 socks*") {
    throw "Edit 1: paragraph 40 did not contain the expected socks sentence: $($socksRange.Text)"
}

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/>' +
  '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/></w:rPr>' +
  '<w:t xml:space="preserve">There are 20 socks in a drawer: 5 pairs of black socks, 3 pairs of brown and 2 pairs of white. You select the socks in the dark and can check them only after a selection has been made. What is the smallest number of socks you need to select to guarantee getting the </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/></w:rPr><w:t>following:</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>'

$socksRange.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: "girl counts" a/b/c paragraphs - wrap the numeral in each line with
# proofErr gramStart/gramEnd, and reflow so each letter keeps its own
# trailing space run.
# ---------------------------------------------------------------------------
$pA = $d.Paragraphs.Item(74)
$pC = $d.Paragraphs.Item(76)
$girlRange = $d.Range($pA.Range.Start, $pC.Range.End)
if ($pA.Range.Text -notlike "a) What if the girl counts from 1 to 10*") {
    throw "Edit 2: paragraph 74 did not contain the expected text: $($pA.Range.Text)"
}
if ($pC.Range.Text -notlike "*1000*") {
    throw "Edit 2: paragraph 76 did not contain the expected text: $($pC.Range.Text)"
}

$pPrGirl = '<w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/></w:rPr></w:pPr>'
$rPrGirl  = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/></w:rPr>'

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $pPrGirl +
  '<w:r>' + $rPrGirl + '<w:t xml:space="preserve">a) What if the girl counts from 1 to </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r>' + $rPrGirl + '<w:t>10</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r>' + $rPrGirl + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>' +
  '<w:p>' + $pPrGirl +
  '<w:r>' + $rPrGirl + '<w:t xml:space="preserve">b) What if the girl counts from 1 to </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r>' + $rPrGirl + '<w:t>100</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r>' + $rPrGirl + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>' +
  '<w:p>' + $pPrGirl +
  '<w:r>' + $rPrGirl + '<w:t xml:space="preserve">c) What if the girl counts from 1 to </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r>' + $rPrGirl + '<w:t>1000</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>'

$girlRange.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Edit 3: final "solution steps" paragraph - move the "_GoBack" bookmark out
# of the middle of the sentence (after "them ") to the very end of the
# document, append the finished sentence, and add two new explanatory
# paragraphs (plus 4 blank separator paragraphs) before the final bookmark
# paragraph.
# ---------------------------------------------------------------------------
$pSolStart = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pSolEnd = $d.Paragraphs.Item($d.Paragraphs.Count)
$solRange = $d.Range($pSolStart.Range.Start, $pSolEnd.Range.End)
if ($pSolStart.Range.Text -notlike "The solution would be to follow the steps shown below and repeat them in increments of 10 until you reached the desired number.*") {
    throw "Edit 3: solution paragraph did not contain the expected text: $($pSolStart.Range.Text)"
}

$pPrSol = '<w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="220"/><w:tab w:val="left" w:pos="720"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>'
$rPrSol  = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>'

$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $pPrSol +
  '<w:r>' + $rPrSol + '<w:t>The solution would be to follow the steps</w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve"> shown below</w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve"> and repeat </w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve">them </w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t>in increments of 10 until you reached the</w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve"> desired number</w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t>.</w:t></w:r>' +
  '</w:p>' +
  '<w:p/><w:p/><w:p/><w:p/>' +
  '<w:p>' + $pPrSol +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve">Since every 10 </w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve">lands on the </w:t></w:r>' +
  '<w:r>' + $rPrSol + '<w:t xml:space="preserve">first finger all increments of the number will result with this finger being selected. </w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$solRange.InsertXML($xml3)
